$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = "Rv3425"
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = "PPE57 Rv3425 MTCY78.04c"
$ws.Range("D3").Value = "FUNCTION: Plays a key role in regulating innate and adaptive immune responses through human Toll-like receptor 2 (TLR2). Interacts with TLR2, leading to the subsequent activation of the mitogen-activated protein kinase (MAPK) and nuclear factor kappa B (NF-kappa-B) signaling pathways. Induces macrophage activation by augmenting the expression of several cell surface molecules (CD40, CD80, CD86 and MHC class II) and pro-inflammatory cytokines (TNF-alpha, IL-6 and IL-12p40) within macrophages. Also participates in adaptive immunity by directing Th1-polarised immune responses (PubMed:25586105). Stimulates specific humoral and cellular immune responses in tuberculosis (TB) patients (PubMed:17328725, PubMed:19467342, PubMed:23136116). Induces a strong IgG(1) antibody response and an increased Th1/Th2 type immune response in mice (PubMed:18426397). {ECO:0000269|PubMed:17328725, ECO:0000269|PubMed:18426397, ECO:0000269|PubMed:19467342, ECO:0000269|PubMed:23136116, ECO:0000269|PubMed:25586105}."
$ws.Range("E3").Value = 3

# Row 4
$ws.Range("A4").Value = "Rv0683"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "rpsG rps7 Rv0683 MTV040.11"
$ws.Range("D4").Value = "FUNCTION: One of the primary rRNA binding proteins, it binds directly to 16S rRNA where it nucleates assembly of the head domain of the 30S subunit. Is located at the subunit interface close to the decoding center, probably blocks exit of the E-site tRNA. {ECO:0000255|HAMAP-Rule:MF_00480}."
$ws.Range("E4").Value = 3

# Row 5
$ws.Range("A5").Value = "Rv2157c"
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = "murF Rv2157c MTCY270.11"
$ws.Range("D5").Value = "FUNCTION: Involved in cell wall formation. Catalyzes the final step in the synthesis of UDP-N-acetylmuramoyl-pentapeptide, the precursor of murein. {ECO:0000255|HAMAP-Rule:MF_02019}."
$ws.Range("E5").Value = 3

# Row 6
$ws.Range("A6").Value = "Rv0706"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "rplV Rv0706 MTCY210.25"
$ws.Range("D6").Value = "FUNCTION: This protein binds specifically to 23S rRNA; its binding is stimulated by other ribosomal proteins, e.g. L4, L17, and L20. It is important during the early stages of 50S assembly. It makes multiple contacts with different domains of the 23S rRNA in the assembled 50S subunit and ribosome (By similarity). {ECO:0000255|HAMAP-Rule:MF_01331}.; FUNCTION: The globular domain of the protein is located near the polypeptide exit tunnel on the outside of the subunit, while an extended beta-hairpin is found that lines the wall of the exit tunnel in the center of the 70S ribosome. {ECO:0000255|HAMAP-Rule:MF_01331}."
$ws.Range("E6").Value = 3

# Row 7
$ws.Range("A7").Value = "Rv3117"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "cysA1 cysA Rv3117 MTCY164.27; cysA2 Rv0815c MTV043.07c"
$ws.Range("D7").Value = "FUNCTION: May be a sulfotransferase involved in the formation of thiosulfate. {ECO:0000250}."
$ws.Range("E7").Value = 3

# Row 8
$ws.Range("A8").Value = "Rv0796"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = "Rv0796 MTV042.06; Rv1369c MTCY02B12.03c; Rv1756c MTCY28.22c; Rv1764 MTCY28.30; Rv2106 MTCY261.02; Rv2167c MTCY270.01; Rv2279 MTCY339.31c; Rv2355 MTCY98.24; Rv2479c MTV008.35c; Rv2649 MTCY441.18; Rv2814c MTCY16B7.29; Rv3185 MTV014.29; Rv3187 MTV014.31; Rv3326 MTV016.26; Rv3380c MTV004.38c; Rv3475 MTCY13E12.28"
$ws.Range("D8").Value = "FUNCTION: Involved in the transposition of the insertion sequence."
$ws.Range("E8").Value = 3

# Row 9
$ws.Range("A9").Value = "Rv1199c"
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = "Rv1199c MTCI364.11c; Rv2512c MTCY07A7.18c"
$ws.Range("D9").Value = "FUNCTION: Required for the transposition of the insertion element."
$ws.Range("E9").Value = 3

# Row 10
$ws.Range("A10").Value = "Rv0515"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Rv0515"
$ws.Range("E10").Value = 3

# Row 11
$ws.Range("A11").Value = "Rv0094c"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "Rv0094c Rv3467"
$ws.Range("E11").Value = 3

# Row 12
$ws.Range("A12").Value = "Rv3888c"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "Rv3888c"
$ws.Range("E12").Value = 3

# Row 13
$ws.Range("A13").Value = "Rv2315c"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "Rv2315c"
$ws.Range("E13").Value = 3

# Row 14
$ws.Range("A14").Value = "Rv2314c"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "Rv2314c"
$ws.Range("E14").Value = 3

# Row 15
$ws.Range("A15").Value = "Rv1432"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "Rv1432"
$ws.Range("E15").Value = 3

# Row 16
$ws.Range("A16").Value = "Rv2666"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "Rv2666"
$ws.Range("D16").Value = "FUNCTION: Required for the transposition of the insertion element. {ECO:0000256|ARBA:ARBA00002190, ECO:0000256|RuleBase:RU365089}."
$ws.Range("E16").Value = 3

# Row 17
$ws.Range("A17").Value = "Rv1047"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = "Rv1047 Rv3023c Rv3115"
$ws.Range("D17").Value = "FUNCTION: Required for the transposition of the insertion element. {ECO:0000256|ARBA:ARBA00002190, ECO:0000256|RuleBase:RU365089}."
$ws.Range("E17").Value = 3
